$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 247, shifting existing rows 247-324 down to 248-325.
$ws.Rows.Item(247).Insert()

# Populate the newly inserted row 247 with the new record.
$ws.Cells.Item(247, 1).Value = 3
$ws.Cells.Item(247, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(247, 3).Value = "Coquimbo"
$ws.Cells.Item(247, 4).Value = 45229
$ws.Cells.Item(247, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(247, 5).Value = 5
$ws.Cells.Item(247, 6).Value = 100112026
$ws.Cells.Item(247, 7).Value = "Haba"
$ws.Cells.Item(247, 8).Value = "Sin especificar"
$ws.Cells.Item(247, 9).Value = "Primera"
$ws.Cells.Item(247, 10).Value = 40
$ws.Cells.Item(247, 11).Value = 8000
$ws.Cells.Item(247, 12).Value = 8000
$ws.Cells.Item(247, 13).Value = 8000
$ws.Cells.Item(247, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(247, 15).Value = "Provincia de Petorca"
$ws.Cells.Item(247, 16).Value = 320
$ws.Cells.Item(247, 17).Value = 25
$ws.Cells.Item(247, 18).Value = "Hortaliza"
